$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "Rue"
$ws.Range("C1").Value = "Entreprise"
$ws.Range("D1").Value = "Sans"
$ws.Range("E1").Value = "Domicile"
$ws.Range("F1").Value = "Parking"

# Row 2
$ws.Range("B2").Value = "['invgauss', (4.036204719120207, -0.16161763541930152, 0.9612597581681632)]"
$ws.Range("C2").Value = "['beta', (52.595462995358545, 27.703319219454222, -35.25571817939934, 64.73662886568977)]"
$ws.Range("D2").Value = "['beta', (0.5117940197816365, 19.2230902247954, -5.458435128107077e-24, 16.73836142353016)]"
$ws.Range("E2").Value = "['beta', (0.5474258078179828, 0.9000121427405545, -3.7860672030373568e-28, 21.599798314360253)]"
$ws.Range("F2").Value = "['invgamma', (1.2485500680449202, -0.23360330017153055, 1.298576682245404)]"

# Row 3
$ws.Range("B3").Value = "['chi2', (1.3400808769713695, -5.127648667537803e-29, 1.8549957545294835)]"
$ws.Range("C3").Value = "['invgamma', (54.707261432752674, -13.673652829752161, 977.6742062851458)]"
$ws.Range("D3").Value = "['beta', (0.30450496225678947, 88.52544141212059, -3.1546313693862455e-28, 176.8979751292736)]"
$ws.Range("E3").Value = "['beta', (0.547862902931137, 0.6163748511340661, -4.316342930087755e-25, 12.900063624422803)]"
$ws.Range("F3").Value = "['pareto', (3.474502202630147, -5.752794160706152, 5.752794160706151)]"

# Row 4
$ws.Range("B4").Value = "['chi2', (1.3380657374339826, -2.022860775150514e-28, 1.8494210080453803)]"
$ws.Range("C4").Value = "['gamma', (1.1060959455943804, -0.004070015051645493, 2.182753636547731)]"
$ws.Range("D4").Value = "['beta', (0.8611652542892512, 79.35796489004701, -1.5419396633588808e-25, 26.695635778426045)]"
$ws.Range("E4").Value = "['beta', (1.8969335796940112, 1.1440186443623233, -1.0101306135905779, 10.760848245640155)]"
$ws.Range("F4").Value = "['invgamma', (2.030103317853264, -0.273339157562145, 1.9277887930808184)]"

# Row 5
$ws.Range("B5").Value = "['beta', (0.7222192851774487, 1.4655167671346607, -3.876253811786544e-29, 6.985813350784625)]"
$ws.Range("C5").Value = "['beta', (0.8699981002941614, 1.3044457209426563, 0.08333333333330016, 5.986391159320796)]"
$ws.Range("D5").Value = "['beta', (0.3021326639337568, 266.1616492102139, -2.002933488547417e-29, 70.98725221298376)]"
$ws.Range("E5").Value = "['beta', (1.9493353130904008, 1.1828119830488304, -0.7288590762220521, 7.595995794064398)]"
$ws.Range("F5").Value = "['invgauss', (1.3223013343950307, -0.16386888664083987, 1.362835073756798)]"
